$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 1.889766
$ws.Range("H2").Value = 5.669298
$ws.Range("I2").Value = 0.02282384478878274
$ws.Range("J2").Value = 0.02282384478878274
$ws.Range("M2").Value = 4.191702
$ws.Range("N2").Value = 12.575106
$ws.Range("O2").Value = 0.5459604994504406
$ws.Range("P2").Value = 0.5459604994504407
$ws.Range("Q2").Value = 7.921335921732001
$ws.Range("R2").Value = 71.29202329558801
$ws.Range("S2").Value = 0.01246091770026316
$ws.Range("T2").Value = 0.01246091770026316
# Row 3
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 1.889766
$ws.Range("H3").Value = 5.669298
$ws.Range("I3").Value = 0.02282384478878274
$ws.Range("J3").Value = 0.02282384478878274
$ws.Range("M3").Value = 1.788624333333333
$ws.Range("N3").Value = 5.365873
$ws.Range("O3").Value = 0.2329646130273282
$ws.Range("P3").Value = 0.2329646130273283
$ws.Range("Q3").Value = 3.380081451906
$ws.Range("R3").Value = 30.420733067154
$ws.Range("S3").Value = 0.005317148169014574
$ws.Range("T3").Value = 0.005317148169014576
# Row 4
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 1.889766
$ws.Range("H4").Value = 5.669298
$ws.Range("I4").Value = 0.02282384478878274
$ws.Range("J4").Value = 0.02282384478878274
$ws.Range("O4").Value = 0.004608562429307936
$ws.Range("P4").Value = 0.004608562429307937
$ws.Range("Q4").Value = 0.066865590378
$ws.Range("R4").Value = 0.601790313402
$ws.Range("S4").Value = 0.0001051851135859399
$ws.Range("T4").Value = 0.0001051851135859399
# Row 5
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 1.889766
$ws.Range("H5").Value = 5.669298
$ws.Range("I5").Value = 0.02282384478878274
$ws.Range("J5").Value = 0.02282384478878274
$ws.Range("M5").Value = 1.661956
$ws.Range("N5").Value = 4.985868
$ws.Range("O5").Value = 0.2164663250929232
$ws.Range("P5").Value = 0.2164663250929232
$ws.Range("Q5").Value = 3.140707942296
$ws.Range("R5").Value = 28.266371480664
$ws.Range("S5").Value = 0.004940593805919067
$ws.Range("T5").Value = 0.004940593805919068
# Row 6
$ws.Range("G6").Value = 66.98490533333333
$ws.Range("I6").Value = 0.8090171389046609
$ws.Range("J6").Value = 0.8090171389046608
$ws.Range("M6").Value = 4.191702
$ws.Range("N6").Value = 12.575106
$ws.Range("O6").Value = 0.5459604994504406
$ws.Range("P6").Value = 0.5459604994504407
$ws.Range("Q6").Value = 280.780761655544
$ws.Range("R6").Value = 2527.026854899896
$ws.Range("S6").Value = 0.4416914012203551
$ws.Range("T6").Value = 0.4416914012203551
# Row 7
$ws.Range("G7").Value = 66.98490533333333
$ws.Range("I7").Value = 0.8090171389046609
$ws.Range("J7").Value = 0.8090171389046608
$ws.Range("M7").Value = 1.788624333333333
$ws.Range("N7").Value = 5.365873
$ws.Range("O7").Value = 0.2329646130273282
$ws.Range("P7").Value = 0.2329646130273283
$ws.Range("Q7").Value = 119.8108316452298
$ws.Range("R7").Value = 1078.297484807068
$ws.Range("S7").Value = 0.1884723646974006
$ws.Range("T7").Value = 0.1884723646974006
# Row 8
$ws.Range("G8").Value = 66.98490533333333
$ws.Range("I8").Value = 0.8090171389046609
$ws.Range("J8").Value = 0.8090171389046608
$ws.Range("O8").Value = 0.004608562429307936
$ws.Range("P8").Value = 0.004608562429307937
$ws.Range("Q8").Value = 2.370126905409333
$ws.Range("S8").Value = 0.00372840599102222
$ws.Range("T8").Value = 0.00372840599102222
# Row 9
$ws.Range("G9").Value = 66.98490533333333
$ws.Range("I9").Value = 0.8090171389046609
$ws.Range("J9").Value = 0.8090171389046608
$ws.Range("M9").Value = 1.661956
$ws.Range("N9").Value = 4.985868
$ws.Range("O9").Value = 0.2164663250929232
$ws.Range("P9").Value = 0.2164663250929232
$ws.Range("Q9").Value = 111.3259653281653
$ws.Range("R9").Value = 1001.933687953488
$ws.Range("S9").Value = 0.1751249669958829
$ws.Range("T9").Value = 0.1751249669958829
# Row 10
$ws.Range("G10").Value = 1.255973333333333
$ws.Range("H10").Value = 3.76792
$ws.Range("I10").Value = 0.01516914814789243
$ws.Range("J10").Value = 0.01516914814789243
$ws.Range("M10").Value = 4.191702
$ws.Range("N10").Value = 12.575106
$ws.Range("O10").Value = 0.5459604994504406
$ws.Range("P10").Value = 0.5459604994504407
$ws.Range("Q10").Value = 5.26466593328
$ws.Range("R10").Value = 47.38199339952
$ws.Range("S10").Value = 0.008281755699061076
$ws.Range("T10").Value = 0.008281755699061078
# Row 11
$ws.Range("G11").Value = 1.255973333333333
$ws.Range("H11").Value = 3.76792
$ws.Range("I11").Value = 0.01516914814789243
$ws.Range("J11").Value = 0.01516914814789243
$ws.Range("M11").Value = 1.788624333333333
$ws.Range("N11").Value = 5.365873
$ws.Range("O11").Value = 0.2329646130273282
$ws.Range("P11").Value = 0.2329646130273283
$ws.Range("Q11").Value = 2.246464466017777
$ws.Range("R11").Value = 20.21818019416
$ws.Range("S11").Value = 0.003533874728227973
$ws.Range("T11").Value = 0.003533874728227973
# Row 12
$ws.Range("G12").Value = 1.255973333333333
$ws.Range("H12").Value = 3.76792
$ws.Range("I12").Value = 0.01516914814789243
$ws.Range("J12").Value = 0.01516914814789243
$ws.Range("O12").Value = 0.004608562429307936
$ws.Range("P12").Value = 0.004608562429307937
$ws.Range("Q12").Value = 0.04444010445333332
$ws.Range("R12").Value = 0.3999609400799999
$ws.Range("S12").Value = 0.00006990796623898313
$ws.Range("T12").Value = 0.00006990796623898313
# Row 13
$ws.Range("G13").Value = 1.255973333333333
$ws.Range("H13").Value = 3.76792
$ws.Range("I13").Value = 0.01516914814789243
$ws.Range("J13").Value = 0.01516914814789243
$ws.Range("M13").Value = 1.661956
$ws.Range("N13").Value = 4.985868
$ws.Range("O13").Value = 0.2164663250929232
$ws.Range("P13").Value = 0.2164663250929232
$ws.Range("Q13").Value = 2.087372417173333
$ws.Range("R13").Value = 18.78635175456
$ws.Range("S13").Value = 0.003283609754364397
$ws.Range("T13").Value = 0.003283609754364397
# Row 14
$ws.Range("G14").Value = 6.683112333333334
$ws.Range("H14").Value = 20.049337
$ws.Range("I14").Value = 0.08071598208561255
$ws.Range("J14").Value = 0.08071598208561255
$ws.Range("M14").Value = 4.191702
$ws.Range("N14").Value = 12.575106
$ws.Range("O14").Value = 0.5459604994504406
$ws.Range("P14").Value = 0.5459604994504407
$ws.Range("Q14").Value = 28.013615333858
$ws.Range("R14").Value = 252.122538004722
$ws.Range("S14").Value = 0.04406773789309384
$ws.Range("T14").Value = 0.04406773789309385
# Row 15
$ws.Range("G15").Value = 6.683112333333334
$ws.Range("H15").Value = 20.049337
$ws.Range("I15").Value = 0.08071598208561255
$ws.Range("J15").Value = 0.08071598208561255
$ws.Range("M15").Value = 1.788624333333333
$ws.Range("N15").Value = 5.365873
$ws.Range("O15").Value = 0.2329646130273282
$ws.Range("P15").Value = 0.2329646130273283
$ws.Range("Q15").Value = 11.95357734180011
$ws.Range("R15").Value = 107.582196076201
$ws.Range("S15").Value = 0.01880396753169549
$ws.Range("T15").Value = 0.01880396753169549
# Row 16
$ws.Range("G16").Value = 6.683112333333334
$ws.Range("H16").Value = 20.049337
$ws.Range("I16").Value = 0.08071598208561255
$ws.Range("J16").Value = 0.08071598208561255
$ws.Range("O16").Value = 0.004608562429307936
$ws.Range("P16").Value = 0.004608562429307937
$ws.Range("Q16").Value = 0.2364685636903333
$ws.Range("R16").Value = 2.128217073213
$ws.Range("S16").Value = 0.0003719846424844465
$ws.Range("T16").Value = 0.0003719846424844465
# Row 17
$ws.Range("G17").Value = 6.683112333333334
$ws.Range("H17").Value = 20.049337
$ws.Range("I17").Value = 0.08071598208561255
$ws.Range("J17").Value = 0.08071598208561255
$ws.Range("M17").Value = 1.661956
$ws.Range("N17").Value = 4.985868
$ws.Range("O17").Value = 0.2164663250929232
$ws.Range("P17").Value = 0.2164663250929232
$ws.Range("Q17").Value = 11.10703864105733
$ws.Range("R17").Value = 99.963347769516
$ws.Range("S17").Value = 0.01747229201833877
$ws.Range("T17").Value = 0.01747229201833878
# Row 18
$ws.Range("G18").Value = 5.984124666666666
$ws.Range("H18").Value = 17.952374
$ws.Range("I18").Value = 0.07227388607305152
$ws.Range("J18").Value = 0.0722738860730515
$ws.Range("M18").Value = 4.191702
$ws.Range("N18").Value = 12.575106
$ws.Range("O18").Value = 0.5459604994504406
$ws.Range("P18").Value = 0.5459604994504407
$ws.Range("Q18").Value = 25.083667333516
$ws.Range("R18").Value = 225.753006001644
$ws.Range("S18").Value = 0.03945868693766745
$ws.Range("T18").Value = 0.03945868693766745
# Row 19
$ws.Range("G19").Value = 5.984124666666666
$ws.Range("H19").Value = 17.952374
$ws.Range("I19").Value = 0.07227388607305152
$ws.Range("J19").Value = 0.0722738860730515
$ws.Range("M19").Value = 1.788624333333333
$ws.Range("N19").Value = 5.365873
$ws.Range("O19").Value = 0.2329646130273282
$ws.Range("P19").Value = 0.2329646130273283
$ws.Range("Q19").Value = 10.70335099250022
$ws.Range("R19").Value = 96.33015893250199
$ws.Range("S19").Value = 0.01683725790098966
$ws.Range("T19").Value = 0.01683725790098966
# Row 20
$ws.Range("G20").Value = 5.984124666666666
$ws.Range("H20").Value = 17.952374
$ws.Range("I20").Value = 0.07227388607305152
$ws.Range("J20").Value = 0.0722738860730515
$ws.Range("O20").Value = 0.004608562429307936
$ws.Range("P20").Value = 0.004608562429307937
$ws.Range("Q20").Value = 0.2117362830806666
$ws.Range("R20").Value = 1.905626547726
$ws.Range("S20").Value = 0.0003330787159763473
$ws.Range("T20").Value = 0.0003330787159763473
# Row 21
$ws.Range("G21").Value = 5.984124666666666
$ws.Range("H21").Value = 17.952374
$ws.Range("I21").Value = 0.07227388607305152
$ws.Range("J21").Value = 0.0722738860730515
$ws.Range("M21").Value = 1.661956
$ws.Range("N21").Value = 4.985868
$ws.Range("O21").Value = 0.2164663250929232
$ws.Range("P21").Value = 0.2164663250929232
$ws.Range("Q21").Value = 9.945351894514666
$ws.Range("R21").Value = 89.508167050632
$ws.Range("S21").Value = 0.01564486251841806
$ws.Range("T21").Value = 0.01564486251841806
